# Clustering algorithms.xlsx — add PopSize/SLS experiment settings and a
# second sheet for the clustering-algorithm experiment (DBScan notes).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename "Sheet1" -> "GA-innstillinger" and fill in the GA / SLS
#     settings table. -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GA-innstillinger"

# A1 already holds "GA-innstillinger:" from the original workbook - leave it.

# Extra note to the right of the header.
$ws1.Range("E1").Value = "Bitflip-mutering I benchmark. Den andre for feature selection"

$ws1.Range("A2").Value = "PopSize"
$ws1.Range("B2").Value = 150

$ws1.Range("A3").Value = "Elitist niches"
$ws1.Range("B3").Value = 20
$ws1.Range("E3").Value = "Single point crossover"

$ws1.Range("A4").Value = "Tournament size"
$ws1.Range("B4").Value = 2

$ws1.Range("A5").Value = "Pc"
$c = $ws1.Range("B5")
$c.Formula = "=""0.7"""
$c.Copy()
$c.PasteSpecial(-4163)

$ws1.Range("A6").Value = "Pm"
$c = $ws1.Range("B6")
$c.Formula = "=""0.2"""
$c.Copy()
$c.PasteSpecial(-4163)

$ws1.Range("A7").Value = "Init crowding"
$ws1.Range("B7").Value = 0

$ws1.Range("A8").Value = "PID control rate"
$c = $ws1.Range("B8")
$c.Formula = "=""0.02"""
$c.Copy()
$c.PasteSpecial(-4163)

$ws1.Range("A9").Value = "Wanted niches"
$ws1.Range("B9").Value = 25

$ws1.Range("A10").Value = "Max members in niche"
$ws1.Range("B10").Value = 3

$ws1.Range("A12").Value = "SLS-innstillinger:"

# Row 13 used to hold "25 runs per problem instance." - it now moves down to
# row 19, and row 13 becomes the first SLS setting.
$ws1.Range("A13").Value = "Max-flips"
$ws1.Range("B13").Value = 30

$ws1.Range("A14").Value = "Pn"
$c = $ws1.Range("B14")
$c.Formula = "=""0.3"""
$c.Copy()
$c.PasteSpecial(-4163)

$ws1.Range("A15").Value = "Max-flips in greedy SLS"
$ws1.Range("B15").Value = 40

$ws1.Range("A16").Value = "Accepts first improvement"

$ws1.Range("A19").Value = "25 runs per problem instance."

$ws1.Range("A20").Value = "MaxFEs"
$ws1.Range("B20").Value = "100k"

$excel.CutCopyMode = $false

# Leave selection on an empty cell below the table, matching the recorded
# session, then move away to the new sheet.
$ws1.Range("A22").Select()

# --- Sheet 2: new sheet with the clustering-algorithm experiment note. ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Eksperiment 1 - Klyngingsalgori"

$ws2.Range("A1").Value = "DBScan"
$ws2.Range("A1").Font.Bold = $true

$ws2.Range("A2").Value = "Epsilon = 0.02"

$ws2.Range("A2").Select()
